$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Mayo de 2020 a las 22:40"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1787200
$ws.Range("C4").Value = 18739
$ws.Range("D4").Value = 502652
$ws.Range("E4").Value = 1180242
$ws.Range("G4").Value = 976
$ws.Range("H4").Value = 104306

# Row 17 - Canada
$ws.Range("B17").Value = 89390
$ws.Range("C17").Value = 878
$ws.Range("D17").Value = 47443
$ws.Range("E17").Value = 34968

# Row 27 - Ecuador
$ws.Range("B27").Value = 38571
$ws.Range("C27").Value = 100
$ws.Range("D27").Value = 19190
$ws.Range("E27").Value = 16047
$ws.Range("G27").Value = 21
$ws.Range("H27").Value = 3334

# Row 134 - Malta
$ws.Range("E134").Value = 93
$ws.Range("G134").Value = 2
$ws.Range("H134").Value = 9
